$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment Schedule")

# Insert a new blank column before column N ("Late"/"Outstanding" columns shift right).
$ws.Columns("N").Insert()

# Make "Repayment Schedule" the active sheet and select cell R6, matching the
# author's final on-screen state (this also clears tabSelected on the sheet
# that was previously active).
$ws.Activate()
$ws.Range("R6").Select() | Out-Null
